# Backlog.xlsx - "updated the backlog again"
#
# The "Current" sheet backlog table (Table13, A2:D27) had rows 18-27
# (backlog items 22-31) added previously with the "Priority" column (B)
# left blank. Fill those priorities in now.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current")

$ws.Range("B18").Value = 4
$ws.Range("B19").Value = 4
$ws.Range("B20").Value = 3
$ws.Range("B21").Value = 3
$ws.Range("B22").Value = 2
$ws.Range("B23").Value = 2
$ws.Range("B24").Value = 3
$ws.Range("B25").Value = 1
$ws.Range("B26").Value = 2
$ws.Range("B27").Value = 1

# Leave the selection where the user last left it after entering the data.
$ws.Range("B26").Select()
